$d = $word.ActiveDocument

# 1. Split "dužim od punog" into "dužim od " + the new placeholder text,
#    replacing the literal word "punog" with the Jinja placeholder.
$found = $d.Content.Find.Execute("dužim od punog", $true, $false, $false, $false, $false, $true, 1, $false, "dužim od {{ radnog_vremena }}", 2)

# 2. Re-find just the inserted placeholder text and nudge its font so Word
#    splits it into its own run (matching the original run formatting,
#    including the complex-script font) rather than merging it back into
#    the preceding run.
$placeholder = $d.Content
$found2 = $placeholder.Find.Execute("{{ radnog_vremena }}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $placeholder.Font.Name = "Arial"
    $placeholder.Font.NameBi = "Arial"
}
